# Insert a new weekly record at the top of the Rabanito / Vega Modelo de
# Temuco time series: duplicate row 33 (pushing the existing rows 33-70
# down to 34-71) and then overwrite the date (Fecha) and Volumen of the
# newly inserted row with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("33:33").Copy()
$ws.Rows("33:33").Insert()

$ws.Range("D33").Value = 44763
$ws.Range("J33").Value = 65
